$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.839.29'
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").Value = '2.220.08'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("E4").Value = '  -1.74%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '299.15'
$c.ClearFormats()
$ws.Range("E5").Value = '  -2.35%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '90.38'
$c.ClearFormats()
$ws.Range("E6").Value = '  -3.36%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.553'
$c.ClearFormats()
$ws.Range("E7").Value = '  -3.02%  '

$ws.Range("E8").Value = '  -0.58%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.490'
$c.ClearFormats()
$ws.Range("E9").Value = '  -5.60%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '32.91'
$c.ClearFormats()
$ws.Range("E10").Value = '  -4.59%  '

$ws.Range("E11").Value = '  -2.99%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '6.93'
$c.ClearFormats()
$ws.Range("E12").Value = '  -3.35%  '

$ws.Range("E13").Value = '  -0.71%  '

$ws.Range("D14").Value = '2.560.27'
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").Value = '2.219.58'
$ws.Range("E15").Value = '  +0.47%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '13.39'
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.52%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.773'
$c.ClearFormats()
$ws.Range("E17").Value = '  -6.70%  '

$ws.Range("D18").Value = '43.730.77'
$ws.Range("E18").Value = '  -1.06%  '

$ws.Range("D19").Value = '0.0₃0901'
$ws.Range("E19").Value = '  -5.08%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '11.44'
$c.ClearFormats()
$ws.Range("E20").Value = '  -3.81%  '

$ws.Range("E21").Value = '  -6.16%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '64.36'
$c.ClearFormats()
$ws.Range("E22").Value = '  -1.91%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '235.79'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.59%  '

$ws.Range("E24").Value = '  -4.61%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("E26").Value = '  -4.83%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '38.23'
$c.ClearFormats()
$ws.Range("E27").Value = '  +2.22%  '

$ws.Range("E28").Value = '  -0.63%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.34'
$c.ClearFormats()
$ws.Range("E29").Value = '  -3.92%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '153.21'
$c.ClearFormats()
$ws.Range("E30").Value = '  -0.01%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '19.17'
$c.ClearFormats()
$ws.Range("E31").Value = '  -3.50%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.40'
$c.ClearFormats()
$ws.Range("E32").Value = '  -8.30%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0755'
$c.ClearFormats()
$ws.Range("E33").Value = '  -4.67%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.48'
$c.ClearFormats()
$ws.Range("E34").Value = '  -5.59%  '

$ws.Range("E35").Value = '  -1.99%  '

$ws.Range("E36").Value = '  -8.47%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.102'
$c.ClearFormats()
$ws.Range("E37").Value = '  -6.61%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.67'
$c.ClearFormats()
$ws.Range("E38").Value = '  -5.37%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0297'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.29%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.17'
$c.ClearFormats()
$ws.Range("E40").Value = '  -5.74%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.60'
$c.ClearFormats()
$ws.Range("E41").Value = '  -3.76%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '13.14'
$c.ClearFormats()
$ws.Range("E42").Value = '  -8.82%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("D44").Value = '1.832.35'
$ws.Range("E44").Value = '  +3.68%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.ClearFormats()
$ws.Range("E45").Value = '  +13.04%  '

$ws.Range("E46").Value = '  -5.22%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '67.78'
$c.ClearFormats()
$ws.Range("E47").Value = '  -2.14%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '94.13'
$c.ClearFormats()
$ws.Range("E48").Value = '  -4.26%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '72.68'
$c.ClearFormats()
$ws.Range("E49").Value = '  -7.73%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '13.89'
$c.ClearFormats()
$ws.Range("E50").Value = '  -1.29%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.72'
$c.ClearFormats()

